$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: clear the polite_expressions value (was "nan", becomes empty)
$ws.Cells.Item(11, 3).Value = ""

# Row 12 (new): duplicate annotation entry for the same source text as row 9,
# but recorded as a separate "SMY" sentence_purpose annotation needing revision.
$ws.Cells.Item(12, 1).Value = "parisk"
$ws.Cells.Item(12, 2).Value = 3
$ws.Cells.Item(12, 3).Value = "nan"
$ws.Cells.Item(12, 4).Value = "SMY"
$ws.Cells.Item(12, 5).Value = "RES"
$ws.Cells.Item(12, 6).Value = "afe80f3f-3501-40b4-a3d0-1ad1f86c76ec"
$ws.Cells.Item(12, 7).Value = "r1BRfhiab_annotated.xlsx"
$ws.Cells.Item(12, 8).Value = "Not too surprisingly, the standard multiclass losses do not have the desired property, however approaches that reduce multi-class to binary classification at training time do, namely unnormalized models with penalized log Z (self-normalization), the NCE approach, as well as (the natural in the proposed setting) binary classification loss."
$ws.Cells.Item(12, 9).Value = "Needs Revision"
